$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("N2").Value = "2016-12-31 00:00:00"
$ws.Range("O2").Value = 111840683.52
$ws.Range("P2").Value = 321696017.18
$ws.Range("Q2").Value = 207073625.98
$ws.Range("R2").Value = 47.089378303
$ws.Range("S2").Value = 138461253.84
$ws.Range("T2").Value = 138461253.84
$ws.Range("U2").Value = 73.11819913239999
$ws.Range("V2").Value = 5905328.89
$ws.Range("W2").Value = 71367306.84
$ws.Range("X2").Value = -9650310.23
$ws.Range("Y2").Value = 114519865.99
$ws.Range("Z2").Value = 119186483.99
$ws.Range("AA2").Value = 9326816.699999999
$ws.Range("AG2").Value = 796668.5699999999
$ws.Range("AP2").Value = 77.1372174645
$ws.Range("AQ2").Value = 180.747612404895
$ws.Range("AR2").Value = 199.366937072001
$ws.Range("AS2").Value = 107733000.01
$ws.Range("AT2").Value = 167.854217661954
